# Mise a jour de l'application
# Reworks the small "legend" block that used to live in columns D:E (rows 1-14)
# so that it now lives in columns G:H (rows 1, 4-16), with a new "_" separator
# header spanning D1:H1, and three extra blank-but-styled cells on row 23 (D23:F23).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Read the existing "/!\" + explanation text currently stored in columns
#    D and E (rows 2-14) before we start overwriting anything.
# ---------------------------------------------------------------------------
$oldD = @{}
$oldE = @{}
for ($r = 2; $r -le 14; $r++) {
    $oldD[$r] = $ws.Range("D$r").Value2
    $oldE[$r] = $ws.Range("E$r").Value2
}

# ---------------------------------------------------------------------------
# 2. Clear out the old D/E content (rows 1-14). We'll rebuild row 1 and
#    rows 4-16 (columns G/H) further below.
# ---------------------------------------------------------------------------
$ws.Range("D1:E14").ClearContents() | Out-Null

# ---------------------------------------------------------------------------
# 3. New header row: D1:H1 all become the "_" separator, using the same
#    bold/filled style as A1 and C1.
# ---------------------------------------------------------------------------
$ws.Range("C1").Copy() | Out-Null
$ws.Range("D1:H1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("D1:H1").Value = "_"

# ---------------------------------------------------------------------------
# 4. Move the "/!\" warning icon + its explanatory text two rows down and
#    two columns to the right: old row R (2..14) -> new row R+2, columns G/H.
# ---------------------------------------------------------------------------
for ($r = 2; $r -le 14; $r++) {
    $newRow = $r + 2
    if ($null -ne $oldD[$r]) {
        $ws.Range("G$newRow").Value = $oldD[$r]
    }
    if ($null -ne $oldE[$r]) {
        $ws.Range("H$newRow").Value = $oldE[$r]
    }
}

# ---------------------------------------------------------------------------
# 5. Row 23 gains three extra styled-but-empty cells (D23:F23), matching the
#    wrap-text / vertical-top style already used by B23/C23.
# ---------------------------------------------------------------------------
$ws.Range("C23").Copy() | Out-Null
$ws.Range("D23:F23").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("D23:F23").ClearContents() | Out-Null

# ---------------------------------------------------------------------------
# 6. Column widths: D, E are brand new; F changes; I (beyond the used range)
#    is new too, inheriting the old F width.
# ---------------------------------------------------------------------------
$ws.Columns("D").ColumnWidth = 11.02
$ws.Columns("E").ColumnWidth = 9.45
$ws.Columns("F").ColumnWidth = 10.59
$ws.Columns("I").ColumnWidth = 15.59

# ---------------------------------------------------------------------------
# 7. Sheet view: zoom out to 75% and select H3.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 75
$ws.Range("H3").Select() | Out-Null
